$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# A leading apostrophe is concatenated to force Excel's text (quote-prefix)
# interpretation, since several price values (e.g. "1.000", "25.829.41")
# would otherwise be auto-coerced into numbers and lose their original formatting.

$ws.Range('D2').Value = "'" + '25.829.41'
$ws.Range('E2').Value = "'" + '  -0.63%  '
$ws.Range('D3').Value = "'" + '1.631.05'
$ws.Range('E3').Value = "'" + '  -0.40%  '
$ws.Range('D4').Value = "'" + '1.000'
$ws.Range('E4').Value = "'" + '  -0.12%  '
$ws.Range('D5').Value = "'" + '215.29'
$ws.Range('E5').Value = "'" + '  +0.32%  '
$ws.Range('D6').Value = "'" + '0.5115'
$ws.Range('E6').Value = "'" + '  +0.53%  '
$ws.Range('E7').Value = "'" + '  -0.03%  '
$ws.Range('D8').Value = "'" + '0.2571'
$ws.Range('E8').Value = "'" + '  +0.15%  '
$ws.Range('D9').Value = "'" + '0.06351'
$ws.Range('E9').Value = "'" + '  -0.04%  '
$ws.Range('D10').Value = "'" + '19.54'
$ws.Range('E10').Value = "'" + '  -0.53%  '
$ws.Range('D11').Value = "'" + '0.07778'
$ws.Range('E11').Value = "'" + '  +0.21%  '
$ws.Range('D12').Value = "'" + '4.260'
$ws.Range('E12').Value = "'" + '  -0.24%  '
$ws.Range('D13').Value = "'" + '1.635.64'
$ws.Range('E13').Value = "'" + '  -0.17%  '
$ws.Range('D14').Value = "'" + '1.852.15'
$ws.Range('E14').Value = "'" + '  -0.64%  '
$ws.Range('D15').Value = "'" + '0.5545'
$ws.Range('E15').Value = "'" + '  +1.98%  '
$ws.Range('D16').Value = "'" + '63.74'
$ws.Range('E16').Value = "'" + '  -0.40%  '
$ws.Range('D17').Value = "'" + '0.0₅7511'
$ws.Range('E17').Value = "'" + '  -2.65%  '
$ws.Range('D18').Value = "'" + '25.837.87'
$ws.Range('E18').Value = "'" + '  -0.62%  '
$ws.Range('E19').Value = "'" + '  +0.00%  '
$ws.Range('D20').Value = "'" + '4.448'
$ws.Range('E20').Value = "'" + '  +0.61%  '
$ws.Range('D21').Value = "'" + '195.05'
$ws.Range('E21').Value = "'" + '  -1.99%  '
$ws.Range('D22').Value = "'" + '9.813'
$ws.Range('E22').Value = "'" + '  -0.97%  '
$ws.Range('D23').Value = "'" + '6.025'
$ws.Range('E23').Value = "'" + '  -0.31%  '
$ws.Range('E24').Value = "'" + '  -0.14%  '
$ws.Range('D25').Value = "'" + '1.883'
$ws.Range('E25').Value = "'" + '  -0.45%  '
$ws.Range('D26').Value = "'" + '141.32'
$ws.Range('E26').Value = "'" + '  +0.12%  '
$ws.Range('D27').Value = "'" + '0.1256'
$ws.Range('E27').Value = "'" + '  +4.44%  '
$ws.Range('E28').Value = "'" + '  -0.29%  '
$ws.Range('D29').Value = "'" + '6.732'
$ws.Range('E29').Value = "'" + '  -1.45%  '
$ws.Range('D30').Value = "'" + '1.239'
$ws.Range('E30').Value = "'" + '  +0.51%  '
$ws.Range('D31').Value = "'" + '0.04878'
$ws.Range('E31').Value = "'" + '  -0.45%  '
$ws.Range('D32').Value = "'" + '3.265'
$ws.Range('E32').Value = "'" + '  +0.10%  '
$ws.Range('D33').Value = "'" + '3.175'
$ws.Range('E33').Value = "'" + '  +0.04%  '
$ws.Range('D34').Value = "'" + '1.550'
$ws.Range('E34').Value = "'" + '  +1.36%  '
$ws.Range('E35').Value = "'" + '  -0.58%  '
$ws.Range('D36').Value = "'" + '0.8983'
$ws.Range('E36').Value = "'" + '  -1.04%  '
$ws.Range('D37').Value = "'" + '0.5542'
$ws.Range('E37').Value = "'" + '  +1.45%  '
$ws.Range('E38').Value = "'" + '  -1.80%  '
$ws.Range('D39').Value = "'" + '1.118.88'
$ws.Range('D40').Value = "'" + '0.01554'
$ws.Range('E40').Value = "'" + '  -0.42%  '
$ws.Range('E41').Value = "'" + '  -0.07%  '
$ws.Range('D42').Value = "'" + '5.554'
$ws.Range('E42').Value = "'" + '  +2.09%  '
$ws.Range('E43').Value = "'" + '  -1.66%  '
$ws.Range('D44').Value = "'" + '97.43'
$ws.Range('E44').Value = "'" + '  -1.62%  '
$ws.Range('D45').Value = "'" + '1.777.36'
$ws.Range('E45').Value = "'" + '  +0.13%  '
$ws.Range('E46').Value = "'" + '  -7.59%  '
$ws.Range('D47').Value = "'" + '0.4426'
$ws.Range('E47').Value = "'" + '  -2.25%  '
$ws.Range('D48').Value = "'" + '0.9980'
$ws.Range('E48').Value = "'" + '  +0.15%  '
$ws.Range('D49').Value = "'" + '54.75'
$ws.Range('E49').Value = "'" + '  -0.32%  '
$ws.Range('D50').Value = "'" + '0.05128'
$ws.Range('E50').Value = "'" + '  +0.16%  '
$ws.Range('D51').Value = "'" + '7.626'
$ws.Range('E51').Value = "'" + '  +4.07%  '
